$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: sku 196265469963, quantity 1
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "196265469963"
$ws.Range("B3").Value = 1

# Update the selection shown in the saved sheet view
$ws.Range("C4").Select()
